$d = $word.ActiveDocument

$replacements = @(
    @("2025-04-01 Tuesday", "2025-04-02 Wednesday"),
    @("73×55=4015", "50×23=1150"),
    @("32×39=1248", "65×12=780"),
    @("97×75=7275", "21×25=525"),
    @("82×11=902", "15×56=840"),
    @("27×42=1134", "93×54=5022"),
    @("53×81=4293", "44×72=3168"),
    @("37×49=1813", "70×55=3850"),
    @("40×26=1040", "79×28=2212"),
    @("27×58=1566", "96×18=1728"),
    @("62×84=5208", "71×99=7029"),
    @("28×65=1820", "92×19=1748"),
    @("47×65=3055", "45×65=2925"),
    @("34×74=2516", "82×77=6314"),
    @("50×76=3800", "14×31=434"),
    @("56×91=5096", "98×91=8918"),
    @("31×54=1674", "58×93=5394"),
    @("46×53=2438", "42×50=2100"),
    @("87×11=957", "25×17=425"),
    @("38×11=418", "23×87=2001"),
    @("87×35=3045", "98×60=5880"),
    @("80×74=5920", "25×18=450"),
    @("56×13=728", "33×44=1452"),
    @("28×96=2688", "57×49=2793"),
    @("68×86=5848", "63×13=819"),
    @("86×81=6966", "18×81=1458")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done replacing $($replacements.Count) text values"
